# Insert one new daily price-observation row for "Perejil" (Vega Modelo de
# Temuco) right before the current row 453, pushing the existing rows
# 453-511 down to 454-512 (dimension grows from A1:R511 to A1:R512).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 453; Excel shifts rows 453:511 down to 454:512.
$ws.Rows.Item(453).Insert()

# Populate the newly inserted row 453 with the new record.
$ws.Range("A453").Value = 10
$ws.Range("B453").Value = "Vega Modelo de Temuco"
$ws.Range("C453").Value = "La Araucanía"
$ws.Range("D453").Value = 45142
$ws.Range("E453").Value = 9
$ws.Range("F453").Value = 100112044
$ws.Range("G453").Value = "Perejil"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 50
$ws.Range("K453").Value = 4000
$ws.Range("L453").Value = 4000
$ws.Range("M453").Value = 4000
$ws.Range("N453").Value = "$/docena de atados (3 kilos)"
$ws.Range("O453").Value = "Provincia de Cautín"
$ws.Range("P453").Value = 1333
$ws.Range("Q453").Value = 3
$ws.Range("R453").Value = "Hortaliza"
